$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("III")

# New "our b=" comparison block in columns M/N (rows 5-7)
$ws.Range("M5").Value = "chveni b= "
$ws.Range("N5").Value = 17

$ws.Range("M6").Value = "mashin N"
$ws.Range("N6").Formula = "=(17*0.16)/((4*PI()*10^(-7))*0.37*COS(ATAN(1/2)))"

$ws.Range("N7").Value = ":(((((("

# Widen column N so the new label/number are visible
$ws.Columns.Item(14).ColumnWidth = 14.28515625

# Move the active selection to where the user ended up editing
$ws.Range("H9").Select()
